$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data now has a new weekly record for "Maracuyá" at the
# "Vega Modelo de Temuco" market that needs to be inserted right before the
# existing last row (row 85). Insert a new row there, which pushes the old
# row 85 down to row 86 (and grows the sheet dimension to A1:T86).
$ws.Range("A85").EntireRow.Insert()

# Fill the newly inserted row 85 with the new record.
$ws.Range("A85").Value = 10
$ws.Range("B85").Value = "Vega Modelo de Temuco"
$ws.Range("C85").Value = "La Araucanía"
$ws.Range("D85").Value = 45041
$ws.Range("E85").Value = 9
$ws.Range("F85").Value = "Fruta"
$ws.Range("G85").Value = 100108
$ws.Range("H85").Value = "Tropicales y subtropicales"
$ws.Range("I85").Value = 100108003
$ws.Range("J85").Value = "Maracuyá"
$ws.Range("K85").Value = "Sin especificar"
$ws.Range("L85").Value = "Primera"
$ws.Range("M85").Value = 25
$ws.Range("N85").Value = 50000
$ws.Range("O85").Value = 50000
$ws.Range("P85").Value = 50000
$ws.Range("Q85").Value = "$/caja 18 kilos"
$ws.Range("R85").Value = "Región de Arica y Parinacota"
$ws.Range("S85").Value = 2778
$ws.Range("T85").Value = 18
